$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Normal" named style, used to reset any stray per-cell formatting back to
# the workbook default (style index 0) once we're done tweaking a cell's
# number format.
$normalStyle = $ws.Cells.Item(145, 2).Style

$newRows = @(
    @{ Row = 146; A = 45454.2916666667; B = 0;   C = 7.30000019073486; D = 7.30000019073486; E = 7.30000019073486; F = 7.30000019073486; G = "7.30000019073486"; H = "VARV.MI" },
    @{ Row = 147; A = 45455.2916666667; B = 0;   C = 7.30000019073486; D = 7.30000019073486; E = 7.30000019073486; F = 7.30000019073486; G = "7.30000019073486"; H = "VARV.MI" },
    @{ Row = 148; A = 45456.4360300926; B = 400; C = 7.30000019073486; D = 7.30000019073486; E = 7.30000019073486; F = 7.30000019073486; G = "7.30000019073486"; H = "VARV.MI" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Column A: serial date/time value, matching the custom
    # "yyyy-mm-dd hh:mm:ss" number format + font already used by the rest of
    # the date column (style index 1 in the workbook).
    $ws.Cells.Item($row, 1).NumberFormat = "yyyy-mm-dd hh:mm:ss"
    $ws.Cells.Item($row, 1).Font.Name = "Calibri"
    $ws.Cells.Item($row, 1).Value = $r.A

    # Column B: plain numeric volume.
    $ws.Cells.Item($row, 2).Value = $r.B

    # Columns C-F: plain numeric OHLC values.
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F

    # Column G (adj_close): stored as text in this workbook, not a number, so
    # force text formatting before assigning, then restore the default style
    # so no stray formatting is left on the cell (matches existing rows,
    # which have no explicit style on column G).
    $ws.Cells.Item($row, 7).NumberFormat = "@"
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 7).Style = $normalStyle

    # Column H: ticker text.
    $ws.Cells.Item($row, 8).Value = $r.H
}
